# Update "countries & provincias Spain" worksheet:
#  - refresh case counts for several Castilla y Leon provinces
#  - re-sort the data block (rows 4-63) by "Casos totales" (col B) descending,
#    since the refreshed counts change the ranking
#  - bump the "Datos actualizados" timestamp in the title cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. apply the updated case counts, looked up by province name -----------
function Set-ProvinceRow {
    param($City, $Total, $Activos, $Recuperados, $Muertes)
    for ($r = 4; $r -le 63; $r++) {
        $name = $ws.Cells.Item($r, 1).Value()
        if ($name -eq $City) {
            $ws.Cells.Item($r, 2).Value = $Total
            $ws.Cells.Item($r, 3).Value = $Activos
            $ws.Cells.Item($r, 4).Value = $Recuperados
            $ws.Cells.Item($r, 5).Value = $Muertes
            return
        }
    }
}

Set-ProvinceRow "Salamanca"   483 42  403 38
Set-ProvinceRow "Valladolid"  410 24  369 17
Set-ProvinceRow "Burgos"      392 41  327 24
Set-ProvinceRow "Leon"        362 21  317 24
Set-ProvinceRow "Segovia"     271 32  212 27
Set-ProvinceRow "Avila"       201 23  163 15
Set-ProvinceRow "Soria"       179 14  152 13
Set-ProvinceRow "Zamora"      100 11  82  7
Set-ProvinceRow "Palencia"    72  13  59  0

# --- 2. re-sort the data block by Casos totales (column B), descending ------
$data = @()
for ($r = 4; $r -le 63; $r++) {
    $row = @{
        name = $ws.Cells.Item($r, 1).Value()
        b    = $ws.Cells.Item($r, 2).Value()
        c    = $ws.Cells.Item($r, 3).Value()
        d    = $ws.Cells.Item($r, 4).Value()
        e    = $ws.Cells.Item($r, 5).Value()
    }
    $data += ,$row
}

$sorted = $data | Sort-Object -Property b -Descending

for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $i + 4
    $row = $sorted[$i]
    $ws.Cells.Item($r, 1).Value = $row.name
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
}

# --- 3. bump the "updated at" timestamp in the title cell -------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 14:16"
